$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above the old row 4, shifting the remaining task rows down,
# and add the new task description in column C (matching the other task rows).
$ws.Rows.Item(4).Insert()
$ws.Range("C4").Value = "Turn section#about other interests ul into an inline-block list with thumbnail pic next to each li that compliments the li"

# Update the active selection to reflect where the cursor ended up after editing.
$ws.Range("C19").Select() | Out-Null
